$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 2P" ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 0
$ws2.Range("F2").Value = 11
$ws2.Range("G2").Value = 100
$ws2.Range("H2").Value = 8.1

$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 1
$ws2.Range("F3").Value = 16
$ws2.Range("G3").Value = 94.12
$ws2.Range("H3").Value = 7.4

$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 1
$ws2.Range("F4").Value = 24
$ws2.Range("G4").Value = 96
$ws2.Range("H4").Value = 7.1

# --- Sheet "Estadisticos Final" ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Range("H2").Value = 8.6
$ws3.Range("H3").Value = 8.1
$ws3.Range("H4").Value = 8.5
